# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# to match the crypto price snapshot captured in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force numeric-looking Price cells to stay as plain text (matches the
# source data's inline-string representation, preserving formats like
# trailing zeros, e.g. "74.00" or "5.30").
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated Price / Volume(1h) values cell by cell.
$ws.Range("D2").Value = "42.656.91"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "2.296.01"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "316.28"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "103.23"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("D10").Value = "39.57"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").Value = "0.0906"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "8.52"
$ws.Range("E12").Value = "  +2.27%  "
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("E14").Value = "  +3.98%  "
$ws.Range("D15").Value = "15.37"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "2.647.87"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "2.297.86"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "42.646.47"
$ws.Range("D19").Value = "7.55"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").Value = "13.75"
$ws.Range("E21").Value = "  +21.50%  "
$ws.Range("D22").Value = "74.00"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("D23").Value = "3.53"
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("D24").Value = "266.65"
$ws.Range("E24").Value = "  -5.03%  "
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("D27").Value = "10.92"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("E28").Value = "  -2.03%  "
$ws.Range("D29").Value = "22.60"
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("D30").Value = "6.58"
$ws.Range("E30").Value = "  +11.90%  "
$ws.Range("D31").Value = "37.42"
$ws.Range("E31").Value = "  +3.63%  "
$ws.Range("D32").Value = "165.51"
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("D33").Value = "0.0882"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("D34").Value = "0.131"
$ws.Range("E34").Value = "  -3.02%  "
$ws.Range("D35").Value = "2.57"
$ws.Range("E35").Value = "  -4.05%  "
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("D37").Value = "4.57"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("D40").Value = "2.72"
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("E41").Value = "  +9.59%  "
$ws.Range("D42").Value = "70.57"
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("D43").Value = "95.83"
$ws.Range("E43").Value = "  -4.71%  "
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "12.48"
$ws.Range("E46").Value = "  +3.96%  "
$ws.Range("D47").Value = "116.82"
$ws.Range("E47").Value = "  +4.41%  "
$ws.Range("D48").Value = "80.02"
$ws.Range("E48").Value = "  +3.89%  "
$ws.Range("D49").Value = "1.667.70"
$ws.Range("E49").Value = "  +4.37%  "
$ws.Range("D50").Value = "5.30"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "8.87"
$ws.Range("E51").Value = "  -0.63%  "
